$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "56 (49 – 61)" "54 (1)"
Replace-Text "55 (48 – 60)" "52 (1)"
Replace-Text "53 (44 – 60)" "51 (1)"
Replace-Text "56 (48 – 61)" "53 (0)"

Replace-Text "n (unweighted) (%); Median (Q1 – Q3)" "n (unweighted) (%); Mean (SE)"

Replace-Text "122 (112 – 134)" "124 (1)"
Replace-Text "128 (116 – 143)" "132 (2)"
Replace-Text "123 (114 – 135)" "125 (2)"
Replace-Text "123 (113 – 135)" "125 (1)"

Replace-Text "72 (66 – 79)" "72 (1)"
Replace-Text "76 (66 – 86)" "76 (1)"
Replace-Text "73 (65 – 79)" "72 (1)"
Replace-Text "73 (66 – 80)" "73 (1)"

Replace-Text "n (unweighted)/N (unweighted) (%); Median (Q1 – Q3)" "n (unweighted)/N (unweighted) (%); Mean (SE)"
